$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.443.40'
$ws.Range('E2').Value = '  +3.88%  '
$ws.Range('D3').Value = '1.590.24'
$ws.Range('E3').Value = '  +2.54%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.994'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.54%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.08'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.512'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.74%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.993'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.64%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '26.62'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +11.49%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.250'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.15%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0595'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.23%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0909'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.30%  '
$ws.Range('D12').Value = '1.812.14'
$ws.Range('E12').Value = '  +2.21%  '
$ws.Range('D13').Value = '1.585.70'
$ws.Range('E13').Value = '  +2.30%  '
$ws.Range('D14').Value = '29.395.68'
$ws.Range('E14').Value = '  +3.74%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.75'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.73%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.528'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.55%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.20'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.99%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '240.80'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +5.88%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.53'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.34%  '
$ws.Range('D20').Value = '0.0₃0692'
$ws.Range('E20').Value = '  +2.61%  '
$ws.Range('E21').Value = '  -0.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.03'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.95%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.29'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.09%  '
$ws.Range('E24').Value = '  +3.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.49'
$ws.Range('D25').Style = 'Normal'
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.23'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.47%  '
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.108'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +5.25%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.39'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.52%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.994'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.55%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0473'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.25%  '
$ws.Range('E31').Value = '  +0.09%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.24'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.58%  '
$ws.Range('D33').Value = '1.428.74'
$ws.Range('E33').Value = '  +3.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.10'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.01%  '
$ws.Range('E35').Value = '  -3.23%  '
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.52'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.87%  '
$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.83'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +9.86%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.31'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.40%  '
$ws.Range('E39').Value = '  +2.72%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.534'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.71%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.98'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.99%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '53.81'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +27.87%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.804'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.35%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.995'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.41%  '
$ws.Range('E45').Value = '  +3.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '64.89'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.81%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.38'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.32%  '
$ws.Range('D48').Value = '1.723.33'
$ws.Range('E48').Value = '  +2.24%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '86.27'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.17%  '
$ws.Range('E50').Value = '  -3.20%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0519'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.87%  '
